# Rename the "Requested quantity" headers on the existing sheets.
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used by the other sheets in the workbook.
$wsForecast.PageSetup.LeftMargin = $wsWeekly.PageSetup.LeftMargin
$wsForecast.PageSetup.RightMargin = $wsWeekly.PageSetup.RightMargin
$wsForecast.PageSetup.TopMargin = $wsWeekly.PageSetup.TopMargin
$wsForecast.PageSetup.BottomMargin = $wsWeekly.PageSetup.BottomMargin
$wsForecast.PageSetup.HeaderMargin = $wsWeekly.PageSetup.HeaderMargin
$wsForecast.PageSetup.FooterMargin = $wsWeekly.PageSetup.FooterMargin

# Header row.
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# Copy the bold/centered header style from an existing sheet's header row.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (date, PO_Forecast, yhat_lower, yhat_upper).
$poForecastData = @(
    @(2, 45221.99999999999, 12, -51.25743914976111, 74.84945949667423),
    @(3, 45228.99999999999, 15, -49.17850245682017, 79.61616259034845),
    @(4, 45235.99999999999, 19, -41.13949474977947, 83.1884317155405),
    @(5, 45242.99999999999, 22, -40.39006702180588, 81.62674818870795),
    @(6, 45249.99999999999, 26, -36.37979099221545, 88.61641505075227),
    @(7, 45256.99999999999, 30, -37.50284701764343, 97.24734838650484),
    @(8, 45263.99999999999, 33, -35.42501021388281, 95.48199375766403),
    @(9, 45270.99999999999, 37, -26.58237071831785, 104.6775882773263),
    @(10, 45277.99999999999, 40, -27.14292373018342, 105.8915463863712),
    @(11, 45298.99999999999, 51, -11.95032034391331, 116.6961991486043),
    @(12, 45305.99999999999, 55, -9.280106421401367, 119.0467902105745),
    @(13, 45312.99999999999, 58, -4.139033554108234, 123.9211580615265),
    @(14, 45319.99999999999, 62, -3.502139977636888, 128.7213909221734),
    @(15, 45326.99999999999, 66, 2.14681333804679, 124.1145285553032),
    @(16, 45333.99999999999, 69, 2.911490013982237, 131.1764577572008),
    @(17, 45347.99999999999, 76, 11.22482532612706, 140.1278424529657),
    @(18, 45354.99999999999, 80, 16.16933101269774, 150.6376800475585),
    @(19, 45361.99999999999, 83, 18.49306164243663, 144.0043814010525),
    @(20, 45368.99999999999, 87, 26.4688450724131, 151.1786457192764),
    @(21, 45375.99999999999, 91, 30.24958360512926, 157.5783486763225),
    @(22, 45382.99999999999, 94, 29.8744125593049, 155.5043041826632),
    @(23, 45389.99999999999, 98, 29.77614212602116, 161.4804639156925),
    @(24, 45396.99999999999, 101, 37.72627230486385, 164.7345353211838),
    @(25, 45403.99999999999, 105, 41.34721860925497, 170.7454049881684),
    @(26, 45410.99999999999, 109, 44.06891335345532, 171.4899362316622),
    @(27, 45417.99999999999, 112, 47.31501243781745, 175.4199129518484),
    @(28, 45424.99999999999, 116, 54.66428146169666, 180.8851488810431),
    @(29, 45431.99999999999, 119, 53.33624083628821, 184.639700921706),
    @(30, 45438.99999999999, 123, 61.25846614650099, 185.7807423458581),
    @(31, 45445.99999999999, 127, 63.09146937344524, 191.0749146290811),
)

foreach ($row in $poForecastData) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

# Apply the date/time number format (matching column A on the other sheets) to the ds column.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A31").PasteSpecial(-4122)

Write-Host "PO Forecast sheet added and headers updated."

